# CustomerSubGroup ve StockSubGroup backend metotları eklendi.
#
# On "Sheet1" a new lookup row is inserted for the stock sub-group route
# code, and a second label is added next to the existing "Cari Grup(40)"
# entry for the customer sub-group route code. This also shifts the
# "Create/Update/Delete/Getlist/GetById" blocks (both the A:C table and
# the G column list) down by one row.
#
# The active sheet/selection also moves from "Sheet6" to "Sheet1".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Insert a new row at 33 - everything from row 33 down (the
# Create(1)/Update(2)/Delete(3)/Getlist(4)/GetById(5) blocks) shifts down
# by one row as a result.
$ws1.Rows.Item(33).Insert()

# Fill the two new cells this edit introduces. A33 is populated before
# B32 so the new shared strings land in the same order as the target
# workbook ("stok sub gruo(60)" then "Cari sub grup(50)").
$ws1.Range("A33").Value = "stok sub gruo(60)"
$ws1.Range("B32").Value = "Cari sub grup(50)"

# The workbook's active tab moves from "Sheet6" to "Sheet1", with the
# selection landing on E23.
$ws1.Activate()
$ws1.Range("E23").Select()
